$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-06-10 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-11 Wednesday", 2)

$d.Content.Find.Execute("266÷4=66, 2", $true, $false, $false, $false, $false, $true, 1, $false, "217÷8=27, 1", 2)
$d.Content.Find.Execute("709÷5=141, 4", $true, $false, $false, $false, $false, $true, 1, $false, "684÷5=136, 4", 2)
$d.Content.Find.Execute("117÷6=19, 3", $true, $false, $false, $false, $false, $true, 1, $false, "433÷3=144, 1", 2)
$d.Content.Find.Execute("934÷9=103, 7", $true, $false, $false, $false, $false, $true, 1, $false, "780÷8=97, 4", 2)
$d.Content.Find.Execute("666÷9=74, 0", $true, $false, $false, $false, $false, $true, 1, $false, "875÷4=218, 3", 2)

$d.Content.Find.Execute("761÷2=380, 1", $true, $false, $false, $false, $false, $true, 1, $false, "874÷2=437, 0", 2)
$d.Content.Find.Execute("177÷4=44, 1", $true, $false, $false, $false, $false, $true, 1, $false, "877÷4=219, 1", 2)
$d.Content.Find.Execute("689÷7=98, 3", $true, $false, $false, $false, $false, $true, 1, $false, "998÷7=142, 4", 2)
$d.Content.Find.Execute("934÷2=467, 0", $true, $false, $false, $false, $false, $true, 1, $false, "662÷8=82, 6", 2)
$d.Content.Find.Execute("922÷4=230, 2", $true, $false, $false, $false, $false, $true, 1, $false, "233÷9=25, 8", 2)

$d.Content.Find.Execute("710÷3=236, 2", $true, $false, $false, $false, $false, $true, 1, $false, "500÷4=125, 0", 2)
$d.Content.Find.Execute("803÷5=160, 3", $true, $false, $false, $false, $false, $true, 1, $false, "818÷2=409, 0", 2)
$d.Content.Find.Execute("469÷8=58, 5", $true, $false, $false, $false, $false, $true, 1, $false, "694÷5=138, 4", 2)
$d.Content.Find.Execute("898÷3=299, 1", $true, $false, $false, $false, $false, $true, 1, $false, "737÷2=368, 1", 2)
$d.Content.Find.Execute("345÷9=38, 3", $true, $false, $false, $false, $false, $true, 1, $false, "524÷7=74, 6", 2)

$d.Content.Find.Execute("683÷3=227, 2", $true, $false, $false, $false, $false, $true, 1, $false, "612÷5=122, 2", 2)
$d.Content.Find.Execute("320÷4=80, 0", $true, $false, $false, $false, $false, $true, 1, $false, "856÷6=142, 4", 2)
$d.Content.Find.Execute("400÷3=133, 1", $true, $false, $false, $false, $false, $true, 1, $false, "236÷8=29, 4", 2)
$d.Content.Find.Execute("300÷5=60, 0", $true, $false, $false, $false, $false, $true, 1, $false, "756÷3=252, 0", 2)
$d.Content.Find.Execute("420÷7=60, 0", $true, $false, $false, $false, $false, $true, 1, $false, "861÷4=215, 1", 2)

$d.Content.Find.Execute("294÷4=73, 2", $true, $false, $false, $false, $false, $true, 1, $false, "142÷6=23, 4", 2)
$d.Content.Find.Execute("936÷2=468, 0", $true, $false, $false, $false, $false, $true, 1, $false, "511÷5=102, 1", 2)
$d.Content.Find.Execute("620÷8=77, 4", $true, $false, $false, $false, $false, $true, 1, $false, "250÷7=35, 5", 2)
$d.Content.Find.Execute("726÷5=145, 1", $true, $false, $false, $false, $false, $true, 1, $false, "846÷6=141, 0", 2)
$d.Content.Find.Execute("268÷8=33, 4", $true, $false, $false, $false, $false, $true, 1, $false, "102÷5=20, 2", 2)
